# Update cryptocurrency price/volume cells to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'26.446.77"
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'  +6.02%  "
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = "'1.720.12"
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'  +3.36%  "
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = "'  +0.26%  "
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = "'331.54"
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = "'1.002"
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'  +0.23%  "
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = "'0.3706"
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'  +2.04%  "
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'48.27"
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'  +1.09%  "
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'0.3349"
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'  +2.29%  "
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = "'1.182"
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'  +4.02%  "
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'0.07375"
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'  +3.85%  "
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'  +0.38%  "
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'6.372"
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'  +5.04%  "
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'20.01"
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'  +2.05%  "
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'7.005"
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'  +5.72%  "
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'1.718.28"
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'  +3.27%  "
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'  +1.66%  "
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'0.06604"
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'  -0.13%  "
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = "'  +3.67%  "
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'  +0.34%  "
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.Value = "'16.51"
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = "'  +4.21%  "
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = "'6.110"
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = "'  +3.03%  "
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = "'12.76"
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'  +1.08%  "
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.Value = "'26.423.60"
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = "'  +6.12%  "
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = "'2.438"
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'  -0.70%  "
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.Value = "'2.383"
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = "'  -2.42%  "
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = "'1.393"
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.Value = "'  +17.16%  "
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.Value = "'151.97"
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.Value = "'  +1.95%  "
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.Value = "'19.30"
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.Value = "'  +3.34%  "
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.Value = "'1.915.50"
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.Value = "'130.66"
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.Value = "'  +3.97%  "
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.Value = "'4.125"
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.Value = "'  +0.95%  "
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.Value = "'5.934"
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.Value = "'  +3.11%  "
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.Value = "'0.08609"
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.Value = "'  +1.54%  "
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.Value = "'1.696"
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.Value = "'  +2.49%  "
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.Value = "'12.62"
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.Value = "'  +2.70%  "
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.Value = "'5.330"
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.Value = "'  +2.59%  "
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.Value = "'0.02313"
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = "'  +1.53%  "
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = "'0.2152"
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'  +3.48%  "
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = "'0.06172"
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'  +0.36%  "
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = "'8.410"
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'  +1.07%  "
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.Value = "'1.226"
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'  -5.18%  "
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.Value = "'0.6160"
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'  +3.40%  "
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'  +0.27%  "
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.Value = "'14.05"
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'  +3.06%  "
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.Value = "'3.900"
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = "'  +1.74%  "
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'0.5941"
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.Value = "'127.60"
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'  +1.76%  "
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'2.027"
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'  +3.60%  "
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'0.07165"
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'  +2.42%  "
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.Value = "'76.58"
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'  +1.96%  "
$c.Style = 'Normal'
